# Update the date line at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-12-20 Friday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-12-21 Saturday", 2)

# Multiplication sign used by all the table entries.
$x = [char]215

# Update each multiplication-problem cell in the table by explicit
# (row, column) position so the known duplicate value ("480x5=2400")
# can't cause an incorrect double-replacement via text search.
$t = $d.Tables.Item(1)

$rows = @(1, 5, 10, 15, 20)
$newValues = @(
    @("514${x}4=2056", "566${x}6=3396", "164${x}6=984",  "778${x}9=7002", "830${x}2=1660"),
    @("275${x}9=2475", "413${x}2=826",  "479${x}7=3353", "721${x}9=6489", "275${x}6=1650"),
    @("533${x}7=3731", "785${x}4=3140", "333${x}2=666",  "348${x}4=1392", "319${x}9=2871"),
    @("277${x}7=1939", "430${x}8=3440", "480${x}5=2400", "217${x}6=1302", "664${x}9=5976"),
    @("389${x}7=2723", "547${x}8=4376", "378${x}9=3402", "616${x}3=1848", "699${x}3=2097")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowIndex = $rows[$i]
    for ($col = 1; $col -le 5; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $newValues[$i][$col - 1]
    }
}
